$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# Remove the obsolete ConverterPath asset row; rows below shift up.
$ws.Rows.Item(4).Delete()

# Make sure every cell in the new block carries the same "empty asset row"
# style (s=2) that already lives on this block before we overwrite the text,
# so cells that previously had no backing <c> element (B9, A11, B11) pick up
# the right formatting too.
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B11").PasteSpecial(4122)

# Add the new asset entries for Treasury sheet id config, master currency
# code & master transaction type, and the service-account key path.
$ws.Range("A8").Value = "SheetIdConfig_Treasury"
$ws.Range("B8").Value = "[Dev] RPA_Moon_SheetIdConfig_Treasury"

$ws.Range("A9").Value = "MasterCurrencyCode"
$ws.Range("B9").Value = "RPA209_VCC_Citibank_MasterCurrencyCode"

$ws.Range("A10").Value = "MasterTransactionType"
$ws.Range("B10").Value = "RPA209_VCC_Citibank_MasterTransactionType"

$ws.Range("A11").Value = "PathSaKey"
$ws.Range("B11").Value = "[Dev] RPA_Moon_PathSaKey"

$ws.Range("A12").Select()
